$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "سرنجات 5 سم" line item (row 13, item #7) was removed from the
# day-sale report. Deleting the whole row shifts the totals row and the
# footer row up by one (old row 14 -> new row 13, old row 15 -> new row 14),
# and Excel also prunes the now-unused shared strings for that item
# ("سرنجات 5 سم", "3.00", "3.0000").
$ws.Rows("13:13").Delete()

# The totals cell (merged P:Q) drops from 269 to 266 to reflect the
# removed item.
$ws.Range("P13").Value = 266

# Deleting a row makes the row that slides into its place inherit the
# height of the row that used to be there; restore row 13's original
# (item-row) height of 24.75pt now that it holds the totals row.
$ws.Rows("13:13").RowHeight = 24.75
